$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price/volume/name/link data per latest scrape.
# Ambiguous numeric-looking price strings are written with a leading
# apostrophe so Excel keeps them as text (preserves e.g. trailing zeros
# like "40.10" / "0.970" instead of coercing to the Number type).

# Row 2
$ws.Range("D2").Value = "42.432.02"
$ws.Range("E2").Value = "  +1.34%  "
# Row 3
$ws.Range("D3").Value = "2.293.42"
$ws.Range("E3").Value = "  +0.71%  "
# Row 4
$ws.Range("E4").Value = "  -0.15%  "
# Row 5
$ws.Range("D5").Value = "'322.72"
$ws.Range("E5").Value = "  +2.51%  "
# Row 6
$ws.Range("D6").Value = "'104.51"
$ws.Range("E6").Value = "  +2.43%  "
# Row 7
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  -0.25%  "
# Row 8
$ws.Range("E8").Value = "  -0.10%  "
# Row 9
$ws.Range("E9").Value = "  +1.36%  "
# Row 10
$ws.Range("D10").Value = "'40.10"
$ws.Range("E10").Value = "  +3.77%  "
# Row 11
$ws.Range("D11").Value = "'0.0909"
$ws.Range("E11").Value = "  +0.97%  "
# Row 12
$ws.Range("D12").Value = "'8.55"
$ws.Range("E12").Value = "  +3.85%  "
# Row 13
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +0.88%  "
# Row 14
$ws.Range("D14").Value = "'0.970"
$ws.Range("E14").Value = "  +1.50%  "
# Row 15
$ws.Range("D15").Value = "'15.28"
$ws.Range("E15").Value = "  +0.76%  "
# Row 16
$ws.Range("D16").Value = "2.638.51"
$ws.Range("E16").Value = "  +0.67%  "
# Row 17
$ws.Range("D17").Value = "2.283.29"
$ws.Range("E17").Value = "  -0.03%  "
# Row 18
$ws.Range("D18").Value = "42.523.24"
$ws.Range("E18").Value = "  +1.68%  "
# Row 19
$ws.Range("D19").Value = "'7.48"
$ws.Range("E19").Value = "  -0.34%  "
# Row 20
$ws.Range("D20").Value = "'0.0000105"
$ws.Range("E20").Value = "  +0.67%  "
# Row 21
$ws.Range("D21").Value = "'13.32"
$ws.Range("E21").Value = "  +35.51%  "
# Row 22
$ws.Range("D22").Value = "'73.29"
$ws.Range("E22").Value = "  +0.07%  "
# Row 23
$ws.Range("E23").Value = "  +1.80%  "
# Row 24
$ws.Range("D24").Value = "'269.59"
$ws.Range("E24").Value = "  -4.89%  "
# Row 25
$ws.Range("D25").Value = "'2.23"
$ws.Range("E25").Value = "  -1.06%  "
# Row 26
$ws.Range("E26").Value = "  -0.44%  "
# Row 27
$ws.Range("D27").Value = "'10.88"
$ws.Range("E27").Value = "  +1.74%  "
# Row 28
$ws.Range("D28").Value = "'2.32"
$ws.Range("E28").Value = "  +1.48%  "
# Row 29
$ws.Range("D29").Value = "'22.56"
$ws.Range("E29").Value = "  -1.44%  "
# Row 30
$ws.Range("D30").Value = "'38.23"
$ws.Range("E30").Value = "  +11.46%  "
# Row 31
$ws.Range("D31").Value = "'165.51"
$ws.Range("E31").Value = "  +1.75%  "
# Row 32
$ws.Range("D32").Value = "'6.18"
$ws.Range("E32").Value = "  +6.87%  "
# Row 33
$ws.Range("D33").Value = "'0.0883"
$ws.Range("E33").Value = "  +1.27%  "
# Row 34
$ws.Range("E34").Value = "  +0.51%  "
# Row 35
$ws.Range("E35").Value = "  -13.13%  "
# Row 36
$ws.Range("E36").Value = "  -0.53%  "
# Row 37
$ws.Range("E37").Value = "  +1.70%  "
# Row 38
$ws.Range("D38").Value = "'0.0355"
$ws.Range("E38").Value = "  +2.97%  "
# Row 39
$ws.Range("D39").Value = "'3.74"
$ws.Range("E39").Value = "  +4.72%  "
# Row 40
$ws.Range("E40").Value = "  -5.60%  "
# Row 41
$ws.Range("E41").Value = "  +6.75%  "
# Row 42
$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D42").Value = "'96.53"
$ws.Range("E42").Value = "  -6.98%  "
# Row 43
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'69.87"
$ws.Range("E43").Value = "  +1.14%  "
# Row 44
$ws.Range("E44").Value = "  -0.17%  "
# Row 45
$ws.Range("E45").Value = "  +0.86%  "
# Row 46
$ws.Range("E46").Value = "  +4.85%  "
# Row 47
$ws.Range("D47").Value = "'81.59"
$ws.Range("E47").Value = "  +8.51%  "
# Row 48
$ws.Range("D48").Value = "'113.41"
$ws.Range("E48").Value = "  -0.99%  "
# Row 49
$ws.Range("E49").Value = "  -0.43%  "
# Row 50
$ws.Range("D50").Value = "'5.27"
$ws.Range("E50").Value = "  +0.00%  "
# Row 51
$ws.Range("D51").Value = "1.593.19"
$ws.Range("E51").Value = "  +4.04%  "
